$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 373; existing rows 373:402 shift down to 374:403.
$ws.Rows.Item(373).Insert()

# Populate the newly inserted row 373 with the new weekly data point.
$ws.Cells.Item(373, 1).Value = 6
$ws.Cells.Item(373, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(373, 3).Value = "Metropolitana"
$ws.Cells.Item(373, 4).Value = 45106
$ws.Cells.Item(373, 5).Value = 13
$ws.Cells.Item(373, 6).Value = 100112026
$ws.Cells.Item(373, 7).Value = "Haba"
$ws.Cells.Item(373, 8).Value = "Sin especificar"
$ws.Cells.Item(373, 9).Value = "Primera"
$ws.Cells.Item(373, 10).Value = 450
$ws.Cells.Item(373, 11).Value = 17000
$ws.Cells.Item(373, 12).Value = 18000
$ws.Cells.Item(373, 13).Value = 17556
$ws.Cells.Item(373, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(373, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(373, 16).Value = 702
$ws.Cells.Item(373, 17).Value = 25
$ws.Cells.Item(373, 18).Value = "Hortaliza"
